$wb = $excel.ActiveWorkbook

# --- workbook.xml: rename the "Organisations" tab to "Organizations" ---
# (sheetPr codeName stays "Organisations" - only the visible tab/sheet name changes)
$wsOrg = $wb.Worksheets.Item("Organisations")
$wsOrg.Name = "Organizations"

# --- sheet2.xml (Files): drop the stray *MISSING-FILE* column ---
$wsFiles = $wb.Worksheets.Item("Files")
$wsFiles.Range("E1").ClearContents()

# --- sheet4.xml (Equipment): "TYPE" -> "TYPE:" ---
$wsEquipment = $wb.Worksheets.Item("Equipment")
$wsEquipment.Range("E1").Value = "TYPE:"

# --- sheet6.xml (Organizations / codeName Organisations) ---
# fix "Organisation" -> "Organization" and add a new row for ANDS
$wsOrg.Range("F2").Value = "Organization"
$wsOrg.Range("A3").Value = "http://ands.org.au"
$wsOrg.Range("B3").Value = "Australian National Data Service"
$wsOrg.Range("C3").Value = "The core purpose iof the Australian National Data Service (ANDS) is to make Australia" + [char]8217 + "s research data assets more valuable for researchers, research institutions and the nation."
$wsOrg.Range("F3").Value = "Organization"
$wsOrg.Hyperlinks.Add($wsOrg.Range("A3"), "http://ands.org.au")

# --- sheet7.xml (Licenses): add the TYPE: column header ---
$wsLicenses = $wb.Worksheets.Item("Licenses")
$wsLicenses.Range("D1").Value = "TYPE:"

# --- sheet8.xml (Publications) ---
$wsPublications = $wb.Worksheets.Item("Publications")
$wsPublications.Range("B1").Value = "Name"
$wsPublications.Range("C1").Value = "RELATION:Creator*"
$wsPublications.Range("E1").Value = "datePublished"
$wsPublications.Range("D2").Value = "ScholarlyArticle"
$wsPublications.Range("E2").Value = 2018

# --- sheet9.xml (Projects) ---
$wsProjects = $wb.Worksheets.Item("Projects")
$wsProjects.Range("D1").Value = "RELATION:Funder*"
$wsProjects.Range("A2").Value = "https://github.com/UTS-eResearch/projects/datacrate"
$wsProjects.Range("B2").Value = "DataCrate Project"
$wsProjects.Range("C2").Value = "The DataCrate project is to write the spec for DataCrate, of which this is an example. The DataCrate project is part of the University of Technology Sydney's Provisioner project."
$wsProjects.Range("A3").Value = "http://eresearch.uts.edu.au/projects/provisioner"
$wsProjects.Range("B3").Value = "Provisioner"
$wsProjects.Range("C3").Value = "The University of Technology Sydney Provisioner project is "
$wsProjects.Range("D3").Value = "University of Technology Sydney, Australian National Data Service"
$wsProjects.Range("E3").Value = "Project "
$wsProjects.Hyperlinks.Add($wsProjects.Range("A2"), "https://github.com/UTS-eResearch/projects/datacrate")
$wsProjects.Hyperlinks.Add($wsProjects.Range("A3"), "http://eresearch.uts.edu.au/projects/provisioner")
